# Add a new "The team" slide at the end of the deck, using the same
# "Title and Content" layout (index 2 in the slide master's layout list)
# as the rest of the content slides (slide8.xml in the target OOXML).

$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# --- Title -------------------------------------------------------------
$s.Shapes.Item(1).TextFrame.TextRange.Text = "The team"

# --- Body / content placeholder -----------------------------------------
# Build the paragraphs with InsertAfter (rather than one big Text= with
# embedded carriage returns) so every run keeps its lang="en-US" marker,
# matching how PowerPoint itself authors new runs.
$body = $s.Shapes.Item(2).TextFrame.TextRange

$body.Text = "Omkar Dokur: User interface"
[void]$body.InsertAfter("`rVyoma Sheth: User interface")
[void]$body.InsertAfter("`rFrancesca Prieto: Literature research and visualizations")
[void]$body.InsertAfter("`rXiangyun (Sherry) Liao: Effect size calculation")
[void]$body.InsertAfter("`rSamira Jahangiri: Reference data wrangler")
[void]$body.InsertAfter("`rJenna Oberstaller: Conceptualization, supervisor, chief documentor")
[void]$body.InsertAfter("`rJustin Gibbons: Conceptualization, supervisor, calculations ")
[void]$body.InsertAfter("`r")

# Bold the name of each team member (everything before the colon).
$body.Characters(1, 6).Font.Bold = $true     # "Omkar "
$body.Characters(7, 5).Font.Bold = $true     # "Dokur"
$body.Characters(29, 5).Font.Bold = $true    # "Vyoma"
$body.Characters(34, 1).Font.Bold = $true    # " "
$body.Characters(35, 5).Font.Bold = $true    # "Sheth"
$body.Characters(57, 16).Font.Bold = $true   # "Francesca Prieto"
$body.Characters(114, 22).Font.Bold = $true  # "Xiangyun (Sherry) Liao"
$body.Characters(162, 7).Font.Bold = $true   # "Samira "
$body.Characters(169, 9).Font.Bold = $true   # "Jahangiri"
$body.Characters(204, 17).Font.Bold = $true  # "Jenna Oberstaller"
$body.Characters(271, 14).Font.Bold = $true  # "Justin Gibbons"
